$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1297.5385
$ws.Range("I41").Value = 875.44446
$ws.Range("J41").Value = 2247.25
$ws.Range("K41").Value = 875.44446
$ws.Range("L41").Value = 2247.25
$ws.Range("M41").Value = -435.44446
$ws.Range("N41").Value = -3127.25
$ws.Range("H112").Value = 3248057.2
$ws.Range("J112").Value = 3497804
$ws.Range("L112").Value = 10493412
$ws.Range("N112").Value = -10495628
$ws.Range("H132").Value = 1161.8975
$ws.Range("I132").Value = 1062.3784
$ws.Range("J132").Value = 3003
$ws.Range("K132").Value = 3187.1352
$ws.Range("L132").Value = 9009
$ws.Range("M132").Value = -657.1352000000002
$ws.Range("N132").Value = -14069
$ws.Range("H138").Value = 1781.9375
$ws.Range("I138").Value = 1468.5358
$ws.Range("K138").Value = 4405.607400000001
$ws.Range("M138").Value = 734.3925999999992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26672.297
$ws.Range("I32").Value = 27206.1
$ws.Range("K32").Value = 27206.1
$ws.Range("M32").Value = -26919.1
$ws.Range("H97").Value = 3368045.2
$ws.Range("I97").Value = 4116227.8
$ws.Range("J97").Value = 1224.5
$ws.Range("K97").Value = 4116227.8
$ws.Range("L97").Value = 1224.5
$ws.Range("M97").Value = -4115731.8
$ws.Range("N97").Value = -2216.5
$ws.Range("H102").Value = 25001330
$ws.Range("I102").Value = 1478.7778
$ws.Range("K102").Value = 1478.7778
$ws.Range("M102").Value = 143.2221999999999
$ws.Range("H132").Value = 2820384.5
$ws.Range("I132").Value = 3392211.2
$ws.Range("J132").Value = 8902.333000000001
$ws.Range("K132").Value = 10176633.6
$ws.Range("L132").Value = 26706.999
$ws.Range("M132").Value = -10174103.6
$ws.Range("N132").Value = -31766.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 34414
$ws.Range("J103").Value = 35218.668
$ws.Range("L103").Value = 35218.668
$ws.Range("N103").Value = -37562.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3903.6904
$ws.Range("I31").Value = 1998.0385
$ws.Range("J31").Value = 7000.375
$ws.Range("K31").Value = 1998.0385
$ws.Range("L31").Value = 7000.375
$ws.Range("M31").Value = -1703.0385
$ws.Range("N31").Value = -7590.375
$ws.Range("H34").Value = 3903.6904
$ws.Range("I34").Value = 1998.0385
$ws.Range("J34").Value = 7000.375
$ws.Range("K34").Value = 1998.0385
$ws.Range("L34").Value = 7000.375
$ws.Range("M34").Value = -1796.0385
$ws.Range("N34").Value = -7404.375
$ws.Range("H55").Value = 34877.11
$ws.Range("I55").Value = 5973.75
$ws.Range("K55").Value = 5973.75
$ws.Range("M55").Value = -5658.75
$ws.Range("H107").Value = 90909670
$ws.Range("I107").Value = 100000510
$ws.Range("K107").Value = 100000510
$ws.Range("M107").Value = -99998590

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 200.4
$ws.Range("J12").Value = 235.42857
$ws.Range("L12").Value = 706.28571
$ws.Range("N12").Value = -1052.28571
$ws.Range("H132").Value = 75197.60000000001
$ws.Range("I132").Value = 122667.78
$ws.Range("K132").Value = 1104010.02
$ws.Range("M132").Value = -1101480.02

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4145.533
$ws.Range("I80").Value = 3375.2222
$ws.Range("J80").Value = 5301
$ws.Range("K80").Value = 3375.2222
$ws.Range("L80").Value = 5301
$ws.Range("M80").Value = -2377.2222
$ws.Range("N80").Value = -7297
$ws.Range("H83").Value = 4145.533
$ws.Range("I83").Value = 3375.2222
$ws.Range("J83").Value = 5301
$ws.Range("K83").Value = 16876.111
$ws.Range("L83").Value = 26505
$ws.Range("M83").Value = -11884.111
$ws.Range("N83").Value = -36489
$ws.Range("H102").Value = 2071.2727
$ws.Range("I102").Value = 2087.111
$ws.Range("K102").Value = 2087.111
$ws.Range("M102").Value = -465.1109999999999
$ws.Range("H107").Value = 700
$ws.Range("I107").Value = 1037.7858
$ws.Range("J107").Value = 362.2143
$ws.Range("K107").Value = 1037.7858
$ws.Range("L107").Value = 362.2143
$ws.Range("M107").Value = 882.2141999999999
$ws.Range("N107").Value = -4202.2143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 79584.766
$ws.Range("I7").Value = 113454.78
$ws.Range("K7").Value = 113454.78
$ws.Range("M7").Value = -113342.78
$ws.Range("H22").Value = 3701.9143
$ws.Range("I22").Value = 2550.1765
$ws.Range("J22").Value = 4789.6665
$ws.Range("K22").Value = 2550.1765
$ws.Range("L22").Value = 4789.6665
$ws.Range("M22").Value = -2255.1765
$ws.Range("N22").Value = -5379.6665
$ws.Range("H27").Value = 3701.9143
$ws.Range("I27").Value = 2550.1765
$ws.Range("J27").Value = 4789.6665
$ws.Range("K27").Value = 2550.1765
$ws.Range("L27").Value = 4789.6665
$ws.Range("M27").Value = -2443.1765
$ws.Range("N27").Value = -5003.6665
$ws.Range("H55").Value = 262.75
$ws.Range("I55").Value = 192.76923
$ws.Range("J55").Value = 392.7143
$ws.Range("K55").Value = 192.76923
$ws.Range("L55").Value = 392.7143
$ws.Range("M55").Value = -19.76922999999999
$ws.Range("N55").Value = -738.7143
$ws.Range("H126").Value = 79584.766
$ws.Range("I126").Value = 113454.78
$ws.Range("K126").Value = 340364.34
$ws.Range("M126").Value = -337894.34
$ws.Range("H136").Value = 1964809
$ws.Range("I136").Value = 2669580.5
$ws.Range("J136").Value = 7110.5557
$ws.Range("K136").Value = 8008741.5
$ws.Range("L136").Value = 21331.6671
$ws.Range("M136").Value = -8006191.5
$ws.Range("N136").Value = -26431.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 196666
$ws.Range("I54").Value = 120000
$ws.Range("J54").Value = 349998
$ws.Range("K54").Value = 120000
$ws.Range("L54").Value = 349998
$ws.Range("M54").Value = -119480
$ws.Range("N54").Value = -351038
$ws.Range("H81").Value = 11509.162
$ws.Range("I81").Value = 4945.7856
$ws.Range("J81").Value = 15504.261
$ws.Range("K81").Value = 9891.5712
$ws.Range("L81").Value = 31008.522
$ws.Range("M81").Value = -8830.5712
$ws.Range("N81").Value = -33130.522
$ws.Range("H84").Value = 11509.162
$ws.Range("I84").Value = 4945.7856
$ws.Range("J84").Value = 15504.261
$ws.Range("K84").Value = 49457.856
$ws.Range("L84").Value = 155042.61
$ws.Range("M84").Value = -44153.856
$ws.Range("N84").Value = -165650.61
$ws.Range("H107").Value = 3117.3333
$ws.Range("I107").Value = 1300.3334
$ws.Range("J107").Value = 4934.3335
$ws.Range("K107").Value = 3901.0002
$ws.Range("L107").Value = 14803.0005
$ws.Range("M107").Value = -1981.0002
$ws.Range("N107").Value = -18643.0005
$ws.Range("H122").Value = 2471.162
$ws.Range("I122").Value = 2411.1614
$ws.Range("J122").Value = 2781.1667
$ws.Range("K122").Value = 7233.4842
$ws.Range("L122").Value = 8343.500100000001
$ws.Range("M122").Value = -4783.4842
$ws.Range("N122").Value = -13243.5001
$ws.Range("H136").Value = 3761800.2
$ws.Range("I136").Value = 4927283
$ws.Range("J136").Value = 6354.8887
$ws.Range("K136").Value = 14781849
$ws.Range("L136").Value = 19064.6661
$ws.Range("M136").Value = -14779299
$ws.Range("N136").Value = -24164.6661
